# "addback some models, create others"
# Re-box the header row, bump the apitegromab price input, clear out the
# now-unused indication/MOA filler cells, tag a few rows with a "Q224"
# source marker, and make Main the active/visible sheet.

$wb = $excel.ActiveWorkbook
$wsMain = $wb.Worksheets.Item("Main")
$wsModel = $wb.Worksheets.Item("Model")

# ---------------------------------------------------------------------
# Main sheet content changes
# ---------------------------------------------------------------------

# Price input bump: 24 -> 35 (flows through K4/K7 formulas automatically)
$wsMain.Range("K2").Value = 35

# C3/D3 ("SMA" / "Myostatin") keep their text but lose their border -
# drop down to the no-border/no-style look.
$wsMain.Range("C3:D3").Borders.LineStyle = 0

# E3:G3, and all of C4:G6 were blank filler cells that are removed
# entirely (no value, no style) in the new layout.
$wsMain.Range("E3:G3").Clear()
$wsMain.Range("C4:G4").Clear()
$wsMain.Range("C5:G5").Clear()
$wsMain.Range("C6:G6").Clear()

# New "Q224" source-quarter tags, right aligned, in column L.
$wsMain.Range("L3").Value = "Q224"
$wsMain.Range("L5").Value = "Q224"
$wsMain.Range("L6").Value = "Q224"
$wsMain.Range("L3,L5,L6").HorizontalAlignment = -4152

# ---------------------------------------------------------------------
# Re-box borders for B2:H7
#   - Row 2 becomes its own fully boxed header row (top+bottom all the
#     way across, left edge at B, right edge at H).
#   - Rows 3-6 just carry the left/right sides of the outer box down.
#   - Row 7 closes the outer box with a bottom edge.
# ---------------------------------------------------------------------

$wsMain.Range("B2:H7").Borders.LineStyle = 0

$wsMain.Range("B2:H2").Borders(9).LineStyle = 1
$wsMain.Range("B2:H2").Borders(8).LineStyle = 1
$wsMain.Range("B2").Borders(7).LineStyle = 1
$wsMain.Range("H2").Borders(10).LineStyle = 1

$wsMain.Range("B3:B6").Borders(7).LineStyle = 1
$wsMain.Range("H3:H6").Borders(10).LineStyle = 1

$wsMain.Range("B7:H7").Borders(9).LineStyle = 1
$wsMain.Range("B7").Borders(7).LineStyle = 1
$wsMain.Range("H7").Borders(10).LineStyle = 1

# ---------------------------------------------------------------------
# View / selection state
# ---------------------------------------------------------------------

$wsMain.Range("L3").Select()
$wsMain.Activate()

$wsModel.Select()
$wsModel.Application.ActiveWindow.Zoom = 220
$wsMain.Select()
